$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-18 from 2023-10-22 (45221)
# to 2023-10-25 (45224), matching the automatic update reflected in the diff.
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
